# Refresh the "想去人数" (interest count) figures and one "最低票价" (lowest
# ticket price) entry to the latest scrape, as published by the gh-pages
# data-generation job (output generated at dd351a1).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1247
$ws1.Range("F3").Value  = 2015
$ws1.Range("F5").Value  = 172
$ws1.Range("F6").Value  = 411
$ws1.Range("F7").Value  = 46
$ws1.Range("F8").Value  = 502
$ws1.Range("F11").Value = 156
$ws1.Range("F12").Value = 786
$ws1.Range("F15").Value = 3955
$ws1.Range("F16").Value = 2739
$ws1.Range("F17").Value = 837
$ws1.Range("G17").Value = "'68"
$ws1.Range("G17").ClearFormats()
$ws1.Range("F18").Value = 609
$ws1.Range("F19").Value = 335
$ws1.Range("F20").Value = 714
$ws1.Range("F21").Value = 1313
$ws1.Range("F23").Value = 677
$ws1.Range("F24").Value = 282

# --- Sheet "本地生活" (Local life) -------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 111

# --- Sheet "全部类型" (All types) ---------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 111
$ws4.Range("F5").Value  = 1247
$ws4.Range("F6").Value  = 2015
$ws4.Range("F8").Value  = 172
$ws4.Range("F9").Value  = 411
$ws4.Range("F10").Value = 46
$ws4.Range("F11").Value = 502
$ws4.Range("F14").Value = 156
$ws4.Range("F15").Value = 786
$ws4.Range("F20").Value = 3955
$ws4.Range("F21").Value = 2739
$ws4.Range("F22").Value = 837
$ws4.Range("G22").Value = "'68"
$ws4.Range("G22").ClearFormats()
$ws4.Range("F23").Value = 609
$ws4.Range("F24").Value = 335
$ws4.Range("F25").Value = 714
$ws4.Range("F26").Value = 1313
$ws4.Range("F28").Value = 677
$ws4.Range("F29").Value = 282
